$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture the French translation text (currently duplicated in U4) before
# the row/column shuffling destroys the source cell.
$frText = $ws.Range("U4").Value2

# Remove the two rows that are being dropped entirely ("Additional Info" /
# "Applicable"), which shifts the "Cases cannot be assigned..." row up to
# row 2.
$ws.Rows("2:3").Delete()

# Remove the now-redundant trailing "fr" column (was column U).
$ws.Columns("U").Delete()

# Re-order the language columns: column E becomes "fr", column F becomes "de".
$ws.Range("E1").Value = "fr"
$ws.Range("F1").Value = "de"

# Column F (German) previously had no value on this row; the German text was
# stored in column E. Move it over, then populate E2 with the French text
# that used to live in column U.
$ws.Range("F2").Value = $ws.Range("E2").Value2
$ws.Range("E2").Value = $frText
